# Adjustment sensitivity EOL RIR
# Updates the "2020" column header to "2030" and refreshes the sensitivity
# values on each of the four material sheets (Neodymium, Dysprosium,
# Copper, Raw silicon).

$wb = $excel.ActiveWorkbook

function Set-SheetValues($Sheet, $Values) {
    foreach ($addr in $Values.Keys) {
        $Sheet.Range($addr).Value = $Values[$addr]
    }
}

# --- Sheet: Neodymium ---
$wsNeodymium = $wb.Worksheets.Item("Neodymium")
Set-SheetValues $wsNeodymium @{
    "C1" = 2030
    "C2" = 0.000001625819899935208
    "D2" = 0.005885782825730004
    "E2" = 0.006492109438477343
    "B3" = 0.000000000002183980475909259
    "C3" = 0.00007919668242014081
    "D3" = 0.005501068836496137
    "E3" = 0.005498836806399566
    "B4" = 0.00000000000003409259119931335
    "C4" = 0.00007154027686980761
    "D4" = 0.003927863988361095
    "E4" = 0.004603584064156631
    "C5" = 0.000000001587093575472107
    "D5" = 0.0002009765232936477
    "E5" = 0.0003658129441505158
}

# --- Sheet: Dysprosium ---
$wsDysprosium = $wb.Worksheets.Item("Dysprosium")
Set-SheetValues $wsDysprosium @{
    "C1" = 2030
}

# --- Sheet: Copper ---
$wsCopper = $wb.Worksheets.Item("Copper")
Set-SheetValues $wsCopper @{
    "C1" = 2030
    "B2" = 0.000003278497091721097
    "C2" = 0.002541871850645688
    "D2" = 0.5379000242912599
    "E2" = 0.4848590208778999
    "B3" = 0.00002229370101113288
    "C3" = 0.009197341277647557
    "D3" = 0.3863510865451887
    "E3" = 0.3402597778765012
    "B4" = 0.00006612099022439717
    "C4" = 0.002447260422877557
    "D4" = 0.2766023524654511
    "E4" = 0.3027710330745292
    "B5" = 0.00002076994439830034
    "C5" = 0.005407093418076853
    "D5" = 0.5124365740692156
    "E5" = 0.3561425361412697
}

# --- Sheet: Raw silicon ---
$wsRawSilicon = $wb.Worksheets.Item("Raw silicon")
Set-SheetValues $wsRawSilicon @{
    "C1" = 2030
    "B2" = 0.000000496652837099915
    "C2" = 0.00003462735458016471
    "D2" = 0.01385766416475183
    "E2" = 0.0139053799804532
    "B3" = 0.000000530035999530297
    "C3" = 0.0001164057132748931
    "D3" = 0.006708009055267363
    "E3" = 0.00648446477114007
    "B4" = 0.000003397047964529607
    "C4" = 0.00003246107553390979
    "D4" = 0.006270289992800229
    "E4" = 0.007025531601665832
    "B5" = 0.000001823860200208514
    "C5" = 0.00004124207012744105
    "D5" = 0.01240009647204811
    "E5" = 0.009771097180398069
}

